# Reorder "Recorded By" (column G) entries so that "System" is listed
# first among the comma-separated recorders, preserving the relative
# order of the remaining names. Only the three known before-values are
# remapped; everything else (e.g. "System" alone, "admin@admin.com, System",
# "dnasr281@gmail.com" alone, "admin@admin.com, dnasr281@gmail.com") is
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 157) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "backup@backdoor.com, System") {
        $cell.Value = "System, backup@backdoor.com"
    }
    elseif ($val -eq "system, backup@backdoor.com, System") {
        $cell.Value = "System, system, backup@backdoor.com"
    }
}
